$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-randomized ranking data (re_rank scores recomputed; a handful of rows'
# age/race/name got reshuffled along with it within the two realeffort groups).

$ws.Range("D2").Value = "60bfcf5805c5ae12a546f9f3"
$ws.Range("E2").Value = "Giana"
$ws.Range("G2").Value = 7.326165517881154
$ws.Range("D3").Value = "5c5882fc5bfe7600011197cb"
$ws.Range("E3").Value = "Colleen"
$ws.Range("G3").Value = 6.035906762210042
$ws.Range("D4").Value = "60b45e9961dd412bfb6780f8"
$ws.Range("E4").Value = "Jewel"
$ws.Range("G4").Value = 6.027117691378983
$ws.Range("D5").Value = "5e96194b0a9fe909389e9f7b"
$ws.Range("E5").Value = "Tina"
$ws.Range("G5").Value = 5.244593786151905
$ws.Range("D6").Value = "60bd88b8fc436774352f53b9"
$ws.Range("E6").Value = "Annes"
$ws.Range("G6").Value = 5.10415078822819
$ws.Range("D7").Value = "60c0e5899d387663c07eb3a4"
$ws.Range("E7").Value = "Nansi"
$ws.Range("G7").Value = 4.080915110249717
$ws.Range("D8").Value = "6077db0613ce87b4a62a78f9"
$ws.Range("E8").Value = "Lori"
$ws.Range("G8").Value = 1.127328457611293
$ws.Range("D9").Value = "608b14a312c099ac00b721b6"
$ws.Range("E9").Value = "Khushi"
$ws.Range("G9").Value = 1.090887131911884
$ws.Range("D10").Value = "5c0e89c6c323400001e6c4a5"
$ws.Range("E10").Value = "Bri"
$ws.Range("C10").Value = 21
$ws.Range("G10").Value = 0.434968325099591
$ws.Range("H10").Value = "Black or African American"
$ws.Range("D11").Value = "60d5775a99b502eec8cf56b4"
$ws.Range("E11").Value = "Shadaisia"
$ws.Range("C11").Value = 30
$ws.Range("G11").Value = 0.2911261319324809
$ws.Range("D12").Value = "60cb36ee9f58331a33cf5506"
$ws.Range("E12").Value = "Shaniek"
$ws.Range("G12").Value = 0.2519625011376062
$ws.Range("D13").Value = "6036f9b3b1842f8b659b18c7"
$ws.Range("E13").Value = "Kellie"
$ws.Range("C13").Value = 32
$ws.Range("G13").Value = 0.100228771449971
$ws.Range("H13").Value = "White"
$ws.Range("D14").Value = "60b091ed11ccda59e3fc7761"
$ws.Range("E14").Value = "Myles"
$ws.Range("G14").Value = 13.01215981117364
$ws.Range("D15").Value = "601d69a993d94008fb2b25dc"
$ws.Range("E15").Value = "Quinterius"
$ws.Range("G15").Value = 8.480054693822751
$ws.Range("D16").Value = "60c2341fe95d71ee52c043f0"
$ws.Range("E16").Value = "Matthew"
$ws.Range("G16").Value = 7.476982114631664
$ws.Range("D17").Value = "5ff8ad350d084e10f500e48a"
$ws.Range("E17").Value = "Drew"
$ws.Range("G17").Value = 7.283504847855609
$ws.Range("D18").Value = "60bf9943e4e04642d4634ecc"
$ws.Range("E18").Value = "Jamarii"
$ws.Range("G18").Value = 5.288466015748947
$ws.Range("D19").Value = "60db4fde6193c50664c9c478"
$ws.Range("E19").Value = "Edosagbe"
$ws.Range("C19").Value = 22
$ws.Range("G19").Value = 5.276213646677543
$ws.Range("H19").Value = "Black or African American"
$ws.Range("D20").Value = "5dd671942b033b5ec8bc97b4"
$ws.Range("E20").Value = "Juan"
$ws.Range("C20").Value = 26
$ws.Range("G20").Value = 5.00143092683078
$ws.Range("H20").Value = "Hispanic"
$ws.Range("D21").Value = "5e2522d6b734b47915f88275"
$ws.Range("E21").Value = "Corey"
$ws.Range("G21").Value = 4.275591191179309
$ws.Range("D22").Value = "60b322994d0b901954690036"
$ws.Range("E22").Value = "Brennan"
$ws.Range("G22").Value = 4.151918691359521
$ws.Range("D23").Value = "6088fc724afd5c008db33e9d"
$ws.Range("E23").Value = "Masuf"
$ws.Range("G23").Value = 3.329972696343859
$ws.Range("D24").Value = "60b83826821417f8e484a207"
$ws.Range("E24").Value = "Eli"
$ws.Range("C24").Value = 29
$ws.Range("G24").Value = 2.177627561568944
$ws.Range("H24").Value = "White"
$ws.Range("D25").Value = "6097b95056caf5ebb2720002"
$ws.Range("E25").Value = "Damian"
$ws.Range("C25").Value = 50
$ws.Range("G25").Value = 2.103275006361129
$ws.Range("H25").Value = "Black or African American"
